$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(64, 96, 65, 36, 27),
    @(5, 41, 77, 96, 4),
    @(62, 5, 90, 69, 68),
    @(74, 31, 29, 88, 92),
    @(88, 79, 56, 32, 78)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $data[$i].Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $data[$i][$j]
    }
}
